# Extend the full results table (selectednetworks) with additional
# architecture_search, pretraining/processmodel and boreal-sites-prediction rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the previous data region (keep header row 1 untouched) so that
# stale cells from rows that shrink (fewer populated columns) do not linger.
$ws.Range("A2:N30").ClearContents()

$rows = @()
  $rows += ,@(0, "MLP0nP2D0R", "mlp", 0, 2, $null, $null, 0, 10000, 0.502784788608551, 0.7629355788230896, 0.4867032170295715, 0.4867032170295715, "selected")
  $rows += ,@(1, "MLP0nP2D0S", "mlp", 0, 2, $null, $null, 0, 10000, 0.4920576214790344, 0.74643874168396, 0.475677102804184, 0.475677102804184, "selected")
  $rows += ,@(2, "MLP0aP3D0R", "mlp", 0, 3, $null, $null, 0, 10000, 0.6566911935806274, 0.7833381891250609, 0.4924018383026123, 0.4924018383026123, "selected")
  $rows += ,@(3, "MLP0aP3D0S", "mlp", 0, 3, $null, $null, 0, 10000, 0.6142827272415161, 0.7515823245048523, 0.4637706279754639, 0.4637706279754639, "selected")
  $rows += ,@(4, "MLP0aP3D1R", "mlp", 0, 3, $null, $null, 1, 10000, 0.6181010007858276, 0.7336056232452393, 0.4738454222679138, 0.4738454222679138, "selected")
  $rows += ,@(5, "MLP0aP3D1S", "mlp", 0, 3, $null, $null, 1, 10000, 0.6610723733901978, 0.7827129364013672, 0.4912989139556885, 0.4912989139556885, "selected")
  $rows += ,@(6, "MLP0aP3D2R", "mlp", 0, 3, $null, $null, 2, 10000, 0.6720488667488098, 0.7379122972488403, 0.4661184251308441, 0.4661184251308441, "selected")
  $rows += ,@(7, "MLP0aP3D2S", "mlp", 0, 3, $null, $null, 2, 10000, 0.6325932741165161, 0.7514798045158386, 0.4698392748832703, 0.4698392748832703, "selected")
  $rows += ,@(8, "CNN0nP2D0R", "cnn", 0, 2, $null, $null, 0, 10000, 0.6911980509757996, 1.32979416847229, 0.8241565823554993, 0.8241565823554993, "selected")
  $rows += ,@(9, "LSTM0nP2D0R", "lstm", 0, 2, $null, $null, 0, 10000, 0.7621762752532959, 1.278318524360657, 0.7848275899887085, 0.7848275899887085, "selected")
  $rows += ,@(10, "RF0", "rf", 0, $null, $null, $null, $null, $null, 0.4524757705116624, 0.7646869400037543, 0.2675646884543371, 0.471043460239115, "selected")
  $rows += ,@(11, "MLP2aP3D1R", "mlp", 2, 3, $null, $null, 1, 10000, 1.379218339920044, 2.12820816040039, 1.576656341552734, 1.576656341552734, "architecture_search")
  $rows += ,@(12, "MLP2aP3D0R", "mlp", 2, 3, $null, $null, 0, 10000, 1.356328010559082, 1.953316330909729, 1.425268888473511, 1.425268888473511, "architecture_search")
  $rows += ,@(13, "MLP2nP2D0R", "mlp", 2, 2, $null, $null, 0, 10000, 1.218623042106628, 1.931608915328979, 1.40399432182312, 1.40399432182312, "architecture_search")
  $rows += ,@(14, "CNN2nP2D0R", "cnn", 2, 2, $null, $null, 0, 10000, 1.818891525268555, 2.689780950546265, 2.058792591094971, 2.058792591094971, "architecture_search")
  $rows += ,@(15, "LSTM2nP2D0R", "lstm", 2, 2, $null, $null, 0, 10000, 1.861822128295898, 2.64080286026001, 2.015515565872192, 2.015515565872192, "architecture_search")
  $rows += ,@(16, "RF2", "rf", 2, 2, $null, $null, $null, $null, 0.460309834240309, 0.7486888345094803, 0.2712954899951439, 0.4650816831083652, "architecture_search")
  $rows += ,@(17, "MLP7D030P0", "mlp", 7, 3, 30, $null, 2, 50000, 1.369324326515198, 1.58485734462738, 0.6471161842346191, 0.6471161842346191, "pretraining")
  $rows += ,@(18, "MLP7D050P0", "mlp", 7, 3, 50, $null, 2, 50000, 1.405500769615173, 1.511763334274292, 0.6225002408027649, 0.6225002408027649, "pretraining")
  $rows += ,@(19, "MLP7D070P0", "mlp", 7, 3, 70, $null, 2, 60000, 1.392964839935303, 1.474447011947632, 0.5971564054489136, 0.5971564054489136, "pretraining")
  $rows += ,@(20, "MLP7D0100P0", "mlp", 7, 3, 100, $null, 2, 80000, 1.412371277809143, 1.479416847229004, 0.5784268975257874, 0.5784268975257874, "pretraining")
  $rows += ,@(21, "MLP8D030P0", "mlp", 8, 3, 30, $null, 2, 50000, 1.280858159065247, 1.571866273880005, 0.5715242624282837, 0.5715242624282837, "pretraining")
  $rows += ,@(22, "MLP8D050P0", "mlp", 8, 3, 50, $null, 2, 50000, 1.381091356277466, 1.577436804771423, 0.5808306336402893, 0.5808306336402893, "pretraining")
  $rows += ,@(23, "MLP8D070P0", "mlp", 8, 3, 70, $null, 2, 60000, 1.424986124038696, 1.523117899894714, 0.588973879814148, 0.588973879814148, "pretraining")
  $rows += ,@(24, "MLP8D0100P0", "mlp", 8, 3, 100, $null, 2, 80000, 1.425348401069641, 1.526179552078247, 0.5742737054824829, 0.5742737054824829, "pretraining")
  $rows += ,@(25, "preles2008hy", "preles", 0, $null, $null, $null, $null, $null, 1.233824349729843, 1.357822380127588, 0.7645135423199171, 0.8625358425384447, "processmodel")
  $rows += ,@(26, "preles2008bk", "preles", 2, $null, $null, $null, $null, $null, 3.126886912605822, 3.329463433960293, 2.367208924018711, 2.558613191546635, "processmodel")
  $rows += ,@(27, "mlp0nP2D0Rbs", "mlp", 0, 2, $null, $null, $null, $null, $null, 0.7071548871534311, $null, 0.5081766289410362, "borealsitesprediction")
  $rows += ,@(28, "prelesbs", "preles", $null, $null, $null, $null, $null, $null, $null, 0.6688095868566116, $null, 0.4229790817339437, "borealsitesprediction")

$nrows = $rows.Count
$ncols = 14
$arr = New-Object 'object[,]' $nrows,$ncols
for ($i = 0; $i -lt $nrows; $i++) {
  for ($j = 0; $j -lt $ncols; $j++) {
    $arr[$i,$j] = $rows[$i][$j]
  }
}

$ws.Range("A2:N30").Value = $arr

# Apply the header/id-column style (bold, bordered, centered) to column A
# for the newly added rows, matching the existing formatting convention
# used for the id column throughout the table.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A2:A30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

Write-Output "applied $($rows.Count) data rows"
